$d = $word.ActiveDocument

# Locate the bullet paragraph that currently precedes "TECHNICAL SKILLS"
# ("Platform impact: Built redistricting system serving 12,847 analysts
# across 89 organizations") and insert the four new achievement bullets
# right after it, before the TECHNICAL SKILLS heading.

$anchorText = "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains($anchorText)) {
        $anchorIndex = $i
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Anchor paragraph not found"
}

$newBullets = @(
    "• Real-time collaboration at national scale",
    "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%",
    "• Increased voter turnout prediction accuracy from 71% to 87%",
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
)

$insertAfter = $d.Paragraphs.Item($anchorIndex).Range
for ($j = 0; $j -lt $newBullets.Length; $j++) {
    $insertAfter.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($anchorIndex + 1 + $j)
    $newPara.Range.Text = $newBullets[$j]
    $insertAfter = $newPara.Range
}

Write-Host "Inserted $($newBullets.Length) paragraphs after paragraph $anchorIndex"
Write-Host "Total paragraphs now: $($d.Paragraphs.Count)"
